$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new outline block (Syntax/Semantics/Soundness/Completeness)
#    right after "Person should evaluate their conclusion to test its
#    validity." paragraph.
# ---------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "Person should evaluate their conclusion to test its validity.") {
        $anchorPara = $i
        break
    }
}
if ($null -eq $anchorPara) {
    throw "Could not find anchor paragraph for new outline block"
}
$insertBase = $anchorPara
$anchor = $d.Paragraphs($anchorPara).Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 1)
$newp.Range.Text = "Syntax:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 1
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 2)
$newp.Range.Text = "About the expressions itself – words and sentences."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 3)
$newp.Range.Text = "Examples:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 3
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 4)
$newp.Range.Text = "‘Bertrand Russell’ is a proper noun."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 5)
$newp.Range.Text = " ‘likes logic’ is a verb phrase. "
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 6)
$newp.Range.Text = "‘Bertrand Russell likes logic’ is a sentence. "
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 7)
$newp.Range.Text = "Combining a proper noun and a verb phrase in this way makes a sentence."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 8)
$newp.Range.Text = "Semantics:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 1
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 9)
$newp.Range.Text = "About the meanings of the expressions."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 10)
$newp.Range.Text = "Examples:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 3
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 11)
$newp.Range.Text = "‘Bertrand Russell’ refers to a British philosopher. "
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 12)
$newp.Range.Text = "‘Bertrand Russell’ refers to Bertrand Russell. "
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 13)
$newp.Range.Text = "‘likes logic’ expresses a property Russell has. "
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 14)
$newp.Range.Text = "‘Bertrand Russell likes logic’ is true."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 4
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 15)
$newp.Range.Text = "Soundness:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 1
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 16)
$newp.Range.Text = "The property of only being able to prove true things."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 17)
$newp.Range.Text = "Logical system is “sound” if and only if the inference rules of the system admit only valid formulas."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 18)
$newp.Range.Text = "Inference rules do not permit invalid conclusions."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 3
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 19)
$newp.Range.Text = "Completeness:"
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 1
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 20)
$newp.Range.Text = "The property of being able to prove all true things."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 21)
$newp.Range.Text = "Logical system is “complete” if and only if all valid formulas can be derived from the axioms and the inference rules."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 2
$anchor = $newp.Range

$anchor.InsertParagraphAfter()
$newp = $d.Paragraphs($insertBase + 22)
$newp.Range.Text = "No valid formulas that can’t be proved."
$newp.Style = "List Paragraph"
$newp.Range.ListFormat.ListLevelNumber = 3
$anchor = $newp.Range

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the document to the end
#    of the new "No valid formulas that can't be proved." paragraph
#    (collapsed bookmark right after the text, before the pilcrow).
# ---------------------------------------------------------------------
$targetParaIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "No valid formulas that can" + [char]8217 + "t be proved.") {
        $targetParaIdx = $i
        break
    }
}
if ($null -eq $targetParaIdx) {
    throw "Could not find target paragraph for bookmark"
}
$targetPara = $d.Paragraphs($targetParaIdx).Range.Duplicate
[void]$targetPara.MoveEnd(1, -1)        # exclude the paragraph mark
$targetPara.Collapse(0)                 # collapse to the end of the text
$targetPara.InsertAfter([char]1)        # temporary marker char so the
                                         # collapsed position isn't exactly
                                         # at the "end of paragraph" edge
                                         # (that boundary mis-resolves to
                                         # document start in this runtime)
$bmPos = $d.Range($targetPara.Start, $targetPara.Start)
$d.Bookmarks.Add("_GoBack", $bmPos)
$tempCharRange = $d.Range($targetPara.Start, $targetPara.Start + 1)
$tempCharRange.Delete()


# ---------------------------------------------------------------------
# 3) Stamp a <w:lastRenderedPageBreak/> marker right before the run that
#    begins "If term1 is a variable and term2 is any type of term, ..."
#    (this run now starts a new rendered page because of the outline
#    block inserted above it).
# ---------------------------------------------------------------------
$pbParaIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("If") -and $d.Paragraphs($i).Range.Text.Contains("is a variable and")) {
        $pbParaIdx = $i
        break
    }
}
if ($null -eq $pbParaIdx) {
    throw "Could not find paragraph for lastRenderedPageBreak"
}
$pbRange = $d.Paragraphs($pbParaIdx).Range
$insertionPoint = $d.Range($pbRange.Start, $pbRange.Start)
$pbXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertionPoint.InsertXML($pbXml)


# ---------------------------------------------------------------------
# 4) Merge the "unify_with_occurs_check/2." and " (performs occurs
#    check)" runs in the final "Command: ..." paragraph into a single
#    run (they already share identical run formatting).
# ---------------------------------------------------------------------
$cmdParaIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Command:")) {
        $cmdParaIdx = $i
        break
    }
}
if ($null -eq $cmdParaIdx) {
    throw "Could not find Command: paragraph"
}
$cmdRange = $d.Paragraphs($cmdParaIdx).Range

$r1 = $cmdRange.Duplicate
$f1 = $r1.Find
$f1.ClearFormatting()
[void]$f1.Execute("unify_with_occurs_check/2.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start = $r1.Start
$run1End = $r1.End

$r2 = $cmdRange.Duplicate
$f2 = $r2.Find
$f2.ClearFormatting()
[void]$f2.Execute(" (performs occurs check)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run2Start = $r2.Start
$run2End = $r2.End

$mergedText = "unify_with_occurs_check/2. (performs occurs check)"
$oldLen = $run1End - $run1Start
$run1Target = $d.Range($run1Start, $run1End)
$run1Target.Text = $mergedText
$shift = $mergedText.Length - $oldLen
$run2TargetStart = $run2Start + $shift
$run2TargetEnd = $run2End + $shift
$run2Target = $d.Range($run2TargetStart, $run2TargetEnd)
$run2Target.Delete()

